# Apply the "added new block fork heuristic" edit to the game_data sheet.
# This replaces the two sample rows (agent_random vs agent_random) with a
# larger set of agent_minimax match results, and adds a new "Depth" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell for column H ---
$ws.Cells.Item(1, 8).Value = "Depth"

# --- Row data (rows 2-9) ---
$rows = @(
    @{
        A = "agent_minimax"; B = "agent_random"; C = 1; D = "agent_minimax";
        E = "0.4764"; F = "h_block_fork|h_center_control"; H = 5;
        G = "2 1 2 0 0 0`n2 2 1 2 1 0`n2 0 0 0 0 0`n1 1 0 0 0 0`n1 1 0 0 0 0`n2 1 0 0 0 0`n2 1 0 0 0 0"
    },
    @{
        A = "agent_minimax"; B = "agent_random"; C = 1; D = "agent_minimax";
        E = "0.3965"; F = "h_block_fork|h_center_control"; H = 5;
        G = "1 2 0 0 0 0`n0 0 0 0 0 0`n1 0 0 0 0 0`n1 2 2 0 0 0`n1 0 0 0 0 0`n1 2 0 0 0 0`n0 0 0 0 0 0"
    },
    @{
        A = "agent_minimax"; B = "agent_minimax"; C = 2; D = "agent_minimax";
        E = "2.4856"; F = "h_block_fork|h_center_control"; H = 5;
        G = "1 1 2 2 0 0`n2 2 2 1 1 2`n2 1 1 2 1 2`n1 2 1 2 2 1`n2 1 1 2 2 1`n1 1 2 1 2 1`n1 1 2 2 0 0"
    },
    @{
        A = "agent_minimax"; B = "agent_minimax"; C = 2; D = "agent_minimax";
        E = "1.6667"; F = "h_block_fork|h_center_control"; H = 5;
        G = "2 1 2 0 0 0`n2 1 2 2 1 1`n1 2 1 2 2 1`n1 1 2 2 1 2`n1 1 0 0 0 0`n2 2 1 0 0 0`n2 1 1 2 1 2"
    },
    @{
        A = "agent_minimax"; B = "agent_minimax"; C = 2; D = "agent_minimax";
        E = "1.6489"; F = "h_block_fork|h_center_control"; H = 5;
        G = "1 1 1 2 2 0`n1 1 2 0 0 0`n2 2 2 1 2 1`n1 1 2 1 2 2`n1 1 2 1 2 0`n2 0 0 0 0 0`n0 0 0 0 0 0"
    },
    @{
        A = "agent_minimax"; B = "agent_minimax"; C = 2; D = "agent_minimax";
        E = "7.3418"; F = "h_block_fork|h_center_control"; H = 5;
        G = "1 1 1 2 1 0`n2 0 0 0 0 0`n2 2 2 1 1 0`n1 2 2 1 2 2`n1 2 0 0 0 0`n2 2 1 1 2 0`n1 1 1 2 1 2"
    },
    @{
        A = "agent_minimax"; B = "agent_random"; C = 1; D = "agent_minimax";
        E = "0.5217"; F = "h_block_fork|h_center_control"; H = 5;
        G = "0 0 0 0 0 0`n2 0 0 0 0 0`n2 2 0 0 0 0`n1 1 1 1 0 0`n1 0 0 0 0 0`n2 0 0 0 0 0`n0 0 0 0 0 0"
    },
    @{
        A = "agent_minimax"; B = "agent_random"; C = 1; D = "agent_minimax";
        E = "0.2982"; F = "h_block_fork|h_center_control"; H = 5;
        G = "1 0 0 0 0 0`n2 0 0 0 0 0`n0 0 0 0 0 0`n1 1 1 1 0 0`n2 0 0 0 0 0`n2 2 1 2 0 0`n1 2 0 0 0 0"
    }
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    # Match Time is stored as text in the sheet (e.g. "0.4764"), not a number.
    $ws.Cells.Item($r, 5).Value = "'" + $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $r++
}

# --- Update selection to A1 (matches the post-edit sheetView state) ---
[void]$ws.Range("A1").Select()
